# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.
# 展览 sheet: rows 2-6 map to F2..F6
# 全部类型 sheet: rows 2-5 map to F2..F5, and row 7 (F7) corresponds to the
#   same event as 展览!F6 (全部类型 has an extra 演出 row inserted at row 6).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1443
$ws1.Range("F3").Value = 3046
$ws1.Range("F4").Value = 36
$ws1.Range("F5").Value = 487
$ws1.Range("F6").Value = 286

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1443
$ws4.Range("F3").Value = 3046
$ws4.Range("F4").Value = 36
$ws4.Range("F5").Value = 487
$ws4.Range("F7").Value = 286
